$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the computed values in A1:B4 (row 5 is unchanged)
$ws.Range("A1").Value = 0.038589317575456111
$ws.Range("B1").Value = -0.038589318721242466
$ws.Range("A2").Value = -0.043061166389564247
$ws.Range("B2").Value = 0.043061165312743011
$ws.Range("A3").Value = 0.06649506656870699
$ws.Range("B3").Value = -0.066495067661306292
$ws.Range("A4").Value = -0.014596383123959179
$ws.Range("B4").Value = 0.014596381960536024

# Narrow both columns to the same new width (was 14.7109375 / 15.42578125,
# both become 14.42578125). The host quantizes ColumnWidth to 1/6-character
# increments, so 13.6667 is the input that lands closest (14.5 stored).
$ws.Columns.Item(1).ColumnWidth = 13.6667
$ws.Columns.Item(2).ColumnWidth = 13.6667
